# This workbook's data rows got reordered: the full contents of certain
# rows were moved to other row positions (the row numbers/cell formatting
# stay fixed, but which observation record occupies each row changes).
# Below we reproduce that by capturing each affected row's values first,
# then writing them into their new row positions according to the
# permutation cycles observed in the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data columns span A (1) .. AY (51)
$firstCol = 1
$lastCol = 51

function Get-RowValues($row) {
    $vals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals += $ws.Cells.Item($row, $c).Value2()
    }
    return $vals
}

function Set-RowValues($row, $vals) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $new = $vals[$c - $firstCol]
        $old = $ws.Cells.Item($row, $c).Value2()
        # Only touch cells whose value actually changes, to avoid
        # unnecessary re-typing (e.g. date-looking text getting
        # reinterpreted as a date) of untouched data.
        if ($old -ne $new) {
            $ws.Cells.Item($row, $c).Value2 = $new
        }
    }
}

function Rotate-Rows($cycleRows) {
    # $cycleRows is an ordered list of row numbers. The content that is
    # currently in $cycleRows[i] must end up in $cycleRows[i+1] (wrapping
    # around at the end), matching the permutation cycles in the diff.
    $captured = @{}
    foreach ($r in $cycleRows) {
        $captured[$r] = Get-RowValues $r
    }
    $n = $cycleRows.Length
    for ($i = 0; $i -lt $n; $i++) {
        $srcRow = $cycleRows[$i]
        $dstRow = $cycleRows[($i + 1) % $n]
        Set-RowValues $dstRow $captured[$srcRow]
    }
}

# Permutation cycles taken from the target diff (row number -> row number
# that receives its current content):
#   16 -> 17 -> 16
#   18 -> 20 -> 19 -> 21 -> 18
#   27 -> 30 -> 27
#   32 -> 34 -> 32
#   35 -> 37 -> 36 -> 35
#   39 -> 43 -> 39
#   40 -> 41 -> 40
#   53 -> 54 -> 53
Rotate-Rows @(16, 17)
Rotate-Rows @(18, 20, 19, 21)
Rotate-Rows @(27, 30)
Rotate-Rows @(32, 34)
Rotate-Rows @(35, 37, 36)
Rotate-Rows @(39, 43)
Rotate-Rows @(40, 41)
Rotate-Rows @(53, 54)

Write-Output "Row permutation applied."
